$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Report regenerated for handoff: status moves from "In Translation" to
# "Ready for handoff" everywhere it is shown, and the associated
# timestamps advance.

# Overview sheet: zh-cn / de-de status columns (E, F) + "Latest HO Xliff
# Generate Date" (G)
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-06 15:41:18"

# zh-cn detail sheet: Status (C) + Latest Handoff Datetime (H)
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-06 15:40:54"

# de-de detail sheet: Status (C) + Latest Handback DateTime (H)
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-06 15:41:18"

# The Status column text grew longer ("In Translation" -> "Ready for
# handoff"), so Excel's column autosize widens the Status column on all
# three sheets to fit the new text.
$overview.Columns.Item(5).ColumnWidth = 16.333333333333332
$overview.Columns.Item(6).ColumnWidth = 16.333333333333332
$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333332
$dede.Columns.Item(3).ColumnWidth = 16.333333333333332
